# Adds a small "Index" helper column at the very left of the sheet,
# containing a 0-based row counter for each data row, and aligns the
# remaining header labels with the data that was already present.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remember how many rows of data exist before we start shuffling columns.
$lastRow = $ws.UsedRange.Rows.Count
if ($lastRow -lt 2) { $lastRow = 40 }

# Insert a brand new column A; everything that used to live in A..J now
# lives in B..K.
$ws.Columns("A").Insert()

# Header row: label the new index column and tidy up the headers that
# follow it so they once again describe the data beneath them.
$ws.Range("A1").Value = "Index"
$ws.Range("D1").Value = "Altersklasse"
$ws.Range("E1").Value = "Geschlecht"
$ws.Range("F1").Value = "Mannschafts_ID"
$ws.Range("G1").Value = "Widget ID"
$ws.Range("H1").Value = "Code groß"
$ws.Range("I1").Value = "Code klein"

# The trailing two header cells (old "Code groß"/"Code klein", now shifted
# to J1/K1) are no longer needed.
$ws.Range("J1").ClearContents()
$ws.Range("K1").ClearContents()

# Fill the new Index column with a 0-based counter for every data row.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 1).Value = $r - 2
}
